# Auto-generated edit script: updates crypto price/volume data in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.091.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.57%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.655.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.22%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'216.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.69%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5140"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.31%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.25%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2630"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.33%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.17%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -5.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07709"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'Polkadot"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'4.424"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.78%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.650.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.48%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.882.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.23%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -3.21%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -3.11%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'64.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.129.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.31%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -4.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'191.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.74%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.72%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -5.21%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.29%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'139.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.00%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.1217"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.89%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.191"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.12%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'16.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.429"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.06%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05950"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.58%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.268"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.62%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.554"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.43%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.246"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -5.66%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.84%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.9598"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.25%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.14%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D38").Value = "'0.5643"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -8.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01584"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.76%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.938"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.56%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.8529"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.19%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.24%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.006.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -7.99%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'100.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.797.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.30%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -4.10%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'56.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.60%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.42%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.968"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.83%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05170"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.51%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.4181"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.26%  "
$ws.Range("E51").Style = "Normal"
